$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Enter a dynamic-array TRANSPOSE formula at A46 that spills across A46:AQ46,
# reproducing column A (rows 2:44) as a horizontal list of non-terminal symbols.
$ws.Range("A46").Formula2 = "=TRANSPOSE(A2:A44)"

# Adjust the view to match the author's final selection/scroll position.
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A46").Select()
